$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix SCK for SPI(2): was Port2B, now corrected to Port2C
$ws.Range("G8").Value = "C"

# Fix broken link text to D-series pinouts: png -> jpg
$ws.Range("A15").Value = "https://karpova-lab.github.io/pyControl-D-Series-Breakout/_static/pinouts.jpg"

# Update selection to match target state
[void]$ws.Range("D19").Select()
